# "thêm show details DC và sửa date trong tableview "
# Apply the changes described by the commit: add a second prescription row
# (PRESCRIPTION), add the matching detail row to THUOC_TRONG_TOA, rename /
# restock the INSTRUMENT sheet rows with real product names, and leave
# MEDICINE/DISEASE untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. PRESCRIPTION sheet: add row 3 - "Toa Thuoc so 2" prescription, dated
#    12/07/2022 -> 09/07/2022 (same date format/style as the existing row).
# ---------------------------------------------------------------------
$wsPres = $wb.Worksheets.Item("PRESCRIPTION")
$wsPres.Range("A3").Value = 2
$wsPres.Range("B3").Value = "Toa Thuoc so 2"
$wsPres.Range("C3").Value = 44754
$wsPres.Range("D3").Value = 44751
$wsPres.Range("C3").NumberFormat = $wsPres.Range("C2").NumberFormat
$wsPres.Range("D3").NumberFormat = $wsPres.Range("D2").NumberFormat

# ---------------------------------------------------------------------
# 2. THUOC_TRONG_TOA sheet: add row 13 - detail line for the new
#    prescription ("Toa Thuoc so 2" contains Acemol, Hop, dosage "2 lan 1
#    ngay", quantity 2).
# ---------------------------------------------------------------------
$wsDetail = $wb.Worksheets.Item("THUOC_TRONG_TOA")
$wsDetail.Range("A13").Value = 1
$wsDetail.Range("B13").Value = "Acemol"
$wsDetail.Range("C13").Value = "Hop"
$wsDetail.Range("D13").Value = "2 lan 1 ngay"
$wsDetail.Range("E13").Value = 2

# ---------------------------------------------------------------------
# 3. INSTRUMENT sheet: replace the generic instrument names with the real
#    product names/details (this is the "show details DC" part of the
#    commit), update row 6 to the Covid test kit (with its own font), and
#    widen the columns to fit the longer text.
# ---------------------------------------------------------------------
$wsInstr = $wb.Worksheets.Item("INSTRUMENT")

$wsInstr.Range("B2").Value = "Nhiệt kế điện tử Pharmacity"
$wsInstr.Range("B3").Value = "Nhiệt kế hồng ngoại Urgo"
$wsInstr.Range("B4").Value = "Máy đo nồng độ oxy trong máu iMedicare iOM-A8"
$wsInstr.Range("B5").Value = "Băng dán có gạc Urgosterile "

$wsInstr.Range("B6").Value = "Dụng cụ xét nghiệm nhanh Covid-19 Antigen Self-Test Abbott Panbio Test Kit"
$wsInstr.Range("D6").Value = "Hop"
$wsInstr.Range("E6").Value = "Test Covid"

# Give row 6's new "Hop"/"Test Covid" cells their own explicit font (matches
# the new style added to the workbook for this commit).
$rRow6 = $wsInstr.Range("D6:E6")
$rRow6.Font.Name = "Arial"
$rRow6.Font.Size = 10
$rRow6.Font.Color = 0

# Row heights grew slightly (15.75 instead of 12.75) to fit the longer text.
$wsInstr.Rows.Item(2).RowHeight = 15.75
$wsInstr.Rows.Item(3).RowHeight = 15.75
$wsInstr.Rows.Item(4).RowHeight = 15.75
$wsInstr.Rows.Item(5).RowHeight = 15.75
$wsInstr.Rows.Item(6).RowHeight = 15.75

# Widen the Name / Link columns so the long product descriptions fit.
$wsInstr.Columns.Item(2).ColumnWidth = 66.57
$wsInstr.Columns.Item(5).ColumnWidth = 24.14

# Page setup tweak recorded alongside this edit.
$wsInstr.PageSetup.PaperSize = 9
$wsInstr.PageSetup.Orientation = 1

# INSTRUMENT becomes the active sheet/selection (E6) after these edits.
$wsInstr.Activate()
$wsInstr.Range("E6").Select()
